# Adds a "more viz scripts starting point" textbox (r-graphics.org link)
# to the "Basic Bar Charts" slide, and relabels the frequency-matrix table
# header/terms from generic "word/Term#" placeholders to "Variable/Grp#".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# ------------------------------------------------------------------
# 1) New textbox: "https://r-graphics.org/" placed just under the
#    "Rectangles are representative..." band, to the right of the table,
#    above the picture. AddTextbox takes coordinates in points, so the
#    target EMU offsets/extents (x=1506682, y=1673944, cx=4572000,
#    cy=369332) are divided by 914400/72 = 12700 EMU-per-point.
# ------------------------------------------------------------------
$emuPerPt = 12700.0
$tbLeft   = 1506682 / $emuPerPt
$tbTop    = 1673944 / $emuPerPt
$tbWidth  = 4572000 / $emuPerPt
$tbHeight = 369332  / $emuPerPt

$textBox = $s.Shapes.AddTextbox(1, $tbLeft, $tbTop, $tbWidth, $tbHeight)
$textBox.Name = "TextBox 6"

$tf = $textBox.TextFrame
$tf.WordWrap = -1
$tf.AutoSize = 1
$textBox.Fill.Visible = 0

# Set the full text, then re-stamp each sub-range's Text in place; that
# forces the engine to split the paragraph into three separate runs at
# the same boundaries PowerPoint used ("https://r-" | "graphics.org" |
# "/") without introducing any incidental formatting overrides.
$tr = $tf.TextRange
$tr.Text = "https://r-graphics.org/"
$tr.Characters(1, 10).Text = "https://r-"
$tr.Characters(11, 12).Text = "graphics.org"
$tr.Characters(23, 1).Text = "/"

# Move the new textbox right after "Picture 2" (position 5) and before
# "Isosceles Triangle 6" (originally position 6), matching where it was
# inserted in the underlying XML shape tree.
$textBox.ZOrder(1)
for ($i = 0; $i -lt 5; $i++) {
    $textBox.ZOrder(2)
}

# ------------------------------------------------------------------
# 2) Rename the generic placeholder column/terms in the frequency-matrix
#    table to the more descriptive "Variable" / "Grp#" wording. Look the
#    table shape up by name since inserting the textbox above shifted
#    everyone after "Picture 2" down by one Shapes-collection index.
# ------------------------------------------------------------------
$tbl = $s.Shapes.Item("Table 7").Table
$tbl.Cell(1, 1).Shape.TextFrame.TextRange.Text = "Variable"
$tbl.Cell(2, 1).Shape.TextFrame.TextRange.Text = "Grp1"
$tbl.Cell(3, 1).Shape.TextFrame.TextRange.Text = "Grp2"
$tbl.Cell(4, 1).Shape.TextFrame.TextRange.Text = "Grp3"
$tbl.Cell(6, 1).Shape.TextFrame.TextRange.Text = "Grp_n"
